$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (Wins / Losses / Ties) go right after the existing
# "Unnamed: 28" column (AC), i.e. AD1:AF1. Copy the format of an existing
# header cell (bold font + border + centered/top alignment, style index 1)
# onto the new header cells before writing their text so they match the
# rest of the header row exactly.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-42) gets the same team record: 96 wins, 66 losses, 0 ties.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
